# "transaltion fixes, updated db, compare question card started"
#
# The only authored content change in this commit (the rest of the
# upstream diff is Word re-saving the package with a different feature
# set: pruned namespace declarations, recomputed w:lastRenderedPageBreak
# caches, a dropped w16cid:durableId, and a shorter built-in
# w:latentStyles list - none of which correspond to an editable action
# in the Word object model) is in the "EKRAN GLOWNY GRY" section: the
# author deleted a stale line-break + strikethrough "-" + a long note
# about a font-size bug, plus a red "//TEST" marker, right after
# "...arty w 100%, dopracowanie czasu". Word leaves the cursor's
# "_GoBack" bookmark behind at the edit point.

$d = $word.ActiveDocument

# Anchor the start of the doomed text: just after "...dopracowanie czasu"
# (but before the manual line break that precedes the struck-through "-").
$startRng = $d.Content.Duplicate
$startRng.Find.Execute("arty w 100%, dopracowanie czasu", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startRng.End

# Anchor the end of the doomed text: right after the red "//TEST" run,
# i.e. right before the line break that starts "- dopracowanie czasu dla
# wszystkich kart".
$endRng = $d.Content.Duplicate
$endRng.Find.Execute("//TEST", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endRng.End

# Remove the break + struck "-" + bug note + red "//TEST" runs in one go.
$deadRng = $d.Range($startPos, $endPos)
$deadRng.Delete()

# Word drops a "_GoBack" bookmark at the last edit position.
$gobackRng = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $gobackRng) | Out-Null
